$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.481.66"
$ws.Range("E2").Value = "  +9.29%  "
$ws.Range("D3").Value = "1.613.88"
$ws.Range("E3").Value = "  +9.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9913"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3692"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("E8").Value = "  +11.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07093"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.942"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.672"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001090"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.50%  "
$ws.Range("D17").Value = "1.606.24"
$ws.Range("E17").Value = "  +8.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9921"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06782"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +14.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +12.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.056"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.47%  "
$ws.Range("D24").Value = "22.471.61"
$ws.Range("E24").Value = "  +8.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.389"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.560"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +20.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.91%  "
$ws.Range("E28").Value = "  +13.81%  "
$ws.Range("D29").Value = "1.786.93"
$ws.Range("E29").Value = "  +9.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.050"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.185"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +23.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9560"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +17.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08262"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.664"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +15.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.281"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.277"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.625"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +16.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06124"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("E41").Value = "  +8.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2029"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9914"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5945"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.832"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5724"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.09%  "
$ws.Range("E48").Value = "  +7.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.985"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06822"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.08"
$ws.Range("D51").Style = "Normal"
